$d = $word.ActiveDocument

$replacements = @(
    @{old="908×4="; new="982×7="},
    @{old="856×5="; new="136×7="},
    @{old="956×3="; new="388×8="},
    @{old="785×9="; new="156×3="},
    @{old="461×6="; new="601×2="},
    @{old="293×2="; new="262×2="},
    @{old="678×8="; new="536×3="},
    @{old="332×9="; new="870×7="},
    @{old="429×4="; new="250×6="},
    @{old="851×9="; new="267×8="},
    @{old="995×3="; new="361×7="},
    @{old="668×4="; new="877×2="},
    @{old="859×4="; new="201×6="},
    @{old="540×9="; new="741×5="},
    @{old="660×6="; new="222×9="},
    @{old="750×5="; new="885×7="},
    @{old="358×4="; new="298×4="},
    @{old="689×6="; new="539×9="},
    @{old="221×4="; new="428×3="},
    @{old="779×7="; new="429×8="},
    @{old="211×7="; new="315×6="},
    @{old="705×5="; new="612×4="},
    @{old="966×3="; new="309×6="},
    @{old="605×2="; new="337×2="},
    @{old="234×5="; new="820×4="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
